# Refresh crypto price (column D) and 1h volume change (column E) figures
# to the latest scraped values, row by row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = '@'
$ws.Range("D2").Value = '51.521.02'
$ws.Range("D2").NumberFormat = 'General'
$ws.Cells.Item(2, 4).Style = 'Normal'
$ws.Cells.Item(2, 5).Value = '  -0.80%  '

$ws.Range("D3").NumberFormat = '@'
$ws.Range("D3").Value = '2.778.49'
$ws.Range("D3").NumberFormat = 'General'
$ws.Cells.Item(3, 4).Style = 'Normal'
$ws.Cells.Item(3, 5).Value = '  -0.34%  '

$ws.Cells.Item(4, 5).Value = '  -0.04%  '

$ws.Range("D5").NumberFormat = '@'
$ws.Range("D5").Value = '351.90'
$ws.Range("D5").NumberFormat = 'General'
$ws.Cells.Item(5, 4).Style = 'Normal'
$ws.Cells.Item(5, 5).Value = '  -1.73%  '

$ws.Range("D6").NumberFormat = '@'
$ws.Range("D6").Value = '108.41'
$ws.Range("D6").NumberFormat = 'General'
$ws.Cells.Item(6, 4).Style = 'Normal'
$ws.Cells.Item(6, 5).Value = '  -1.05%  '

$ws.Cells.Item(7, 5).Value = '  -1.27%  '

$ws.Range("D9").NumberFormat = '@'
$ws.Range("D9").Value = '0.619'
$ws.Range("D9").NumberFormat = 'General'
$ws.Cells.Item(9, 4).Style = 'Normal'
$ws.Cells.Item(9, 5).Value = '  +5.22%  '

$ws.Range("D10").NumberFormat = '@'
$ws.Range("D10").Value = '39.32'
$ws.Range("D10").NumberFormat = 'General'
$ws.Cells.Item(10, 4).Style = 'Normal'
$ws.Cells.Item(10, 5).Value = '  -1.12%  '

$ws.Cells.Item(11, 5).Value = '  +1.85%  '

$ws.Cells.Item(12, 5).Value = '  -1.39%  '

$ws.Range("D13").NumberFormat = '@'
$ws.Range("D13").Value = '19.83'
$ws.Range("D13").NumberFormat = 'General'
$ws.Cells.Item(13, 4).Style = 'Normal'
$ws.Cells.Item(13, 5).Value = '  +2.13%  '

$ws.Cells.Item(14, 5).Value = '  +3.36%  '

$ws.Range("D15").NumberFormat = '@'
$ws.Range("D15").Value = '3.215.03'
$ws.Range("D15").NumberFormat = 'General'
$ws.Cells.Item(15, 4).Style = 'Normal'
$ws.Cells.Item(15, 5).Value = '  -0.36%  '

$ws.Range("D16").NumberFormat = '@'
$ws.Range("D16").Value = '2.777.36'
$ws.Range("D16").NumberFormat = 'General'
$ws.Cells.Item(16, 4).Style = 'Normal'
$ws.Cells.Item(16, 5).Value = '  +0.07%  '

$ws.Range("D17").NumberFormat = '@'
$ws.Range("D17").Value = '0.924'
$ws.Range("D17").NumberFormat = 'General'
$ws.Cells.Item(17, 4).Style = 'Normal'
$ws.Cells.Item(17, 5).Value = '  -1.17%  '

$ws.Range("D18").NumberFormat = '@'
$ws.Range("D18").Value = '51.508.50'
$ws.Range("D18").NumberFormat = 'General'
$ws.Cells.Item(18, 4).Style = 'Normal'
$ws.Cells.Item(18, 5).Value = '  -0.79%  '

$ws.Range("D19").NumberFormat = '@'
$ws.Range("D19").Value = '7.74'
$ws.Range("D19").NumberFormat = 'General'
$ws.Cells.Item(19, 4).Style = 'Normal'
$ws.Cells.Item(19, 5).Value = '  +3.93%  '

$ws.Cells.Item(20, 5).Value = '  -0.03%  '

$ws.Range("D21").NumberFormat = '@'
$ws.Range("D21").Value = '13.31'
$ws.Range("D21").NumberFormat = 'General'
$ws.Cells.Item(21, 4).Style = 'Normal'
$ws.Cells.Item(21, 5).Value = '  +1.87%  '

$ws.Range("D22").NumberFormat = '@'
$ws.Range("D22").Value = '0.0₃0965'
$ws.Range("D22").NumberFormat = 'General'
$ws.Cells.Item(22, 4).Style = 'Normal'
$ws.Cells.Item(22, 5).Value = '  -0.84%  '

$ws.Range("D23").NumberFormat = '@'
$ws.Range("D23").Value = '70.32'
$ws.Range("D23").NumberFormat = 'General'
$ws.Cells.Item(23, 4).Style = 'Normal'
$ws.Cells.Item(23, 5).Value = '  +0.22%  '

$ws.Range("D24").NumberFormat = '@'
$ws.Range("D24").Value = '266.55'
$ws.Range("D24").NumberFormat = 'General'
$ws.Cells.Item(24, 4).Style = 'Normal'
$ws.Cells.Item(24, 5).Value = '  -1.03%  '

$ws.Cells.Item(25, 5).Value = '  +0.32%  '

$ws.Range("D26").NumberFormat = '@'
$ws.Range("D26").Value = '0.999'
$ws.Range("D26").NumberFormat = 'General'
$ws.Cells.Item(26, 4).Style = 'Normal'
$ws.Cells.Item(26, 5).Value = '  -0.15%  '

$ws.Cells.Item(27, 5).Value = '  -2.25%  '

$ws.Range("D28").NumberFormat = '@'
$ws.Range("D28").Value = '0.164'
$ws.Range("D28").NumberFormat = 'General'
$ws.Cells.Item(28, 4).Style = 'Normal'
$ws.Cells.Item(28, 5).Value = '  +1.19%  '

$ws.Range("D29").NumberFormat = '@'
$ws.Range("D29").Value = '10.27'
$ws.Range("D29").NumberFormat = 'General'
$ws.Cells.Item(29, 4).Style = 'Normal'
$ws.Cells.Item(29, 5).Value = '  +0.49%  '

$ws.Range("D30").NumberFormat = '@'
$ws.Range("D30").Value = '37.09'
$ws.Range("D30").NumberFormat = 'General'
$ws.Cells.Item(30, 4).Style = 'Normal'
$ws.Cells.Item(30, 5).Value = '  +10.05%  '

$ws.Cells.Item(31, 5).Value = '  -2.74%  '

$ws.Range("D32").NumberFormat = '@'
$ws.Range("D32").Value = '6.17'
$ws.Range("D32").NumberFormat = 'General'
$ws.Cells.Item(32, 4).Style = 'Normal'
$ws.Cells.Item(32, 5).Value = '  +7.96%  '

$ws.Range("D33").NumberFormat = '@'
$ws.Range("D33").Value = '52.01'
$ws.Range("D33").NumberFormat = 'General'
$ws.Cells.Item(33, 4).Style = 'Normal'
$ws.Cells.Item(33, 5).Value = '  +0.24%  '

$ws.Range("D34").NumberFormat = '@'
$ws.Range("D34").Value = '0.0446'
$ws.Range("D34").NumberFormat = 'General'
$ws.Cells.Item(34, 4).Style = 'Normal'
$ws.Cells.Item(34, 5).Value = '  -5.02%  '

$ws.Range("D35").NumberFormat = '@'
$ws.Range("D35").Value = '5.55'
$ws.Range("D35").NumberFormat = 'General'
$ws.Cells.Item(35, 4).Style = 'Normal'
$ws.Cells.Item(35, 5).Value = '  +6.91%  '

$ws.Range("D36").NumberFormat = '@'
$ws.Range("D36").Value = '1.00'
$ws.Range("D36").NumberFormat = 'General'
$ws.Cells.Item(36, 4).Style = 'Normal'
$ws.Cells.Item(36, 5).Value = '  -0.04%  '

$ws.Range("D37").NumberFormat = '@'
$ws.Range("D37").Value = '0.0838'
$ws.Range("D37").NumberFormat = 'General'
$ws.Cells.Item(37, 4).Style = 'Normal'
$ws.Cells.Item(37, 5).Value = '  -0.15%  '

$ws.Range("D38").NumberFormat = '@'
$ws.Range("D38").Value = '18.64'
$ws.Range("D38").NumberFormat = 'General'
$ws.Cells.Item(38, 4).Style = 'Normal'
$ws.Cells.Item(38, 5).Value = '  -0.97%  '

$ws.Cells.Item(39, 5).Value = '  -3.34%  '

$ws.Cells.Item(40, 5).Value = '  -1.43%  '

$ws.Cells.Item(41, 5).Value = '  -0.76%  '

$ws.Range("D42").NumberFormat = '@'
$ws.Range("D42").Value = '2.49'
$ws.Range("D42").NumberFormat = 'General'
$ws.Cells.Item(42, 4).Style = 'Normal'
$ws.Cells.Item(42, 5).Value = '  -2.61%  '

$ws.Range("D43").NumberFormat = '@'
$ws.Range("D43").Value = '120.73'
$ws.Range("D43").NumberFormat = 'General'
$ws.Cells.Item(43, 4).Style = 'Normal'
$ws.Cells.Item(43, 5).Value = '  +1.03%  '

$ws.Range("D44").NumberFormat = '@'
$ws.Range("D44").Value = '22.13'
$ws.Range("D44").NumberFormat = 'General'
$ws.Cells.Item(44, 4).Style = 'Normal'
$ws.Cells.Item(44, 5).Value = '  +1.79%  '

$ws.Cells.Item(45, 5).Value = '  -2.29%  '

$ws.Range("D46").NumberFormat = '@'
$ws.Range("D46").Value = '2.146.76'
$ws.Range("D46").NumberFormat = 'General'
$ws.Cells.Item(46, 4).Style = 'Normal'
$ws.Cells.Item(46, 5).Value = '  +3.11%  '

$ws.Range("D47").NumberFormat = '@'
$ws.Range("D47").Value = '3.29'
$ws.Range("D47").NumberFormat = 'General'
$ws.Cells.Item(47, 4).Style = 'Normal'
$ws.Cells.Item(47, 5).Value = '  +1.86%  '

$ws.Cells.Item(48, 5).Value = '  +5.38%  '

$ws.Range("D49").NumberFormat = '@'
$ws.Range("D49").Value = '0.222'
$ws.Range("D49").NumberFormat = 'General'
$ws.Cells.Item(49, 4).Style = 'Normal'
$ws.Cells.Item(49, 5).Value = '  +16.87%  '

$ws.Cells.Item(50, 5).Value = '  -5.65%  '

$ws.Range("D51").NumberFormat = '@'
$ws.Range("D51").Value = '0.899'
$ws.Range("D51").NumberFormat = 'General'
$ws.Cells.Item(51, 4).Style = 'Normal'
$ws.Cells.Item(51, 5).Value = '  -5.52%  '
